# Auto-generated Excel COM-interop script
# Applies market-price/profit recalculations to the Goblin_Profits workbook
# (scheduled runner update) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3582.0557
$ws.Range("J38").Value = 6332.4443
$ws.Range("L38").Value = 18997.3329
$ws.Range("N38").Value = -19741.3329
$ws.Range("H100").Value = 3426.1538
$ws.Range("I100").Value = 1399.0714
$ws.Range("K100").Value = 1399.0714
$ws.Range("M100").Value = -858.0714
$ws.Range("H132").Value = 1730.3784
$ws.Range("I132").Value = 1200.9642
$ws.Range("J132").Value = 3377.4443
$ws.Range("K132").Value = 3602.8926
$ws.Range("L132").Value = 10132.3329
$ws.Range("M132").Value = -1072.8926
$ws.Range("N132").Value = -15192.3329
$ws.Range("H138").Value = 3039.4775
$ws.Range("I138").Value = 2149.8
$ws.Range("K138").Value = 6449.400000000001
$ws.Range("M138").Value = -1309.400000000001
$ws.Range("H141").Value = 5483.222
$ws.Range("I141").Value = 5824.8335
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 17474.5005
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -12294.5005
$ws.Range("N141").Value = -24760
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4782.25
$ws.Range("I32").Value = 4325.4614
$ws.Range("K32").Value = 4325.4614
$ws.Range("M32").Value = -4038.4614
$ws.Range("H61").Value = 4647.45
$ws.Range("I61").Value = 4647.45
$ws.Range("K61").Value = 4647.45
$ws.Range("M61").Value = -4435.45
$ws.Range("H97").Value = 439.45834
$ws.Range("I97").Value = 285.25
$ws.Range("J97").Value = 593.6667
$ws.Range("K97").Value = 285.25
$ws.Range("L97").Value = 593.6667
$ws.Range("M97").Value = 210.75
$ws.Range("N97").Value = -1585.6667
$ws.Range("H102").Value = 3322.1714
$ws.Range("I102").Value = 1867.1786
$ws.Range("K102").Value = 1867.1786
$ws.Range("M102").Value = -245.1786
$ws.Range("H122").Value = 2706.9092
$ws.Range("I122").Value = 2491.5
$ws.Range("K122").Value = 7474.5
$ws.Range("M122").Value = -5024.5
$ws.Range("H136").Value = 4647.45
$ws.Range("I136").Value = 4647.45
$ws.Range("K136").Value = 13942.35
$ws.Range("M136").Value = -11392.35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3002.1667
$ws.Range("I94").Value = 1602.6
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 1602.6
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -1151.6
$ws.Range("N94").Value = -10902
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3168.2424
$ws.Range("I31").Value = 1439.2727
$ws.Range("J31").Value = 6626.1816
$ws.Range("K31").Value = 1439.2727
$ws.Range("L31").Value = 6626.1816
$ws.Range("M31").Value = -1144.2727
$ws.Range("N31").Value = -7216.1816
$ws.Range("H34").Value = 3168.2424
$ws.Range("I34").Value = 1439.2727
$ws.Range("J34").Value = 6626.1816
$ws.Range("K34").Value = 1439.2727
$ws.Range("L34").Value = 6626.1816
$ws.Range("M34").Value = -1237.2727
$ws.Range("N34").Value = -7030.1816
$ws.Range("H47").Value = 16000
$ws.Range("J47").Value = 16000
$ws.Range("L47").Value = 16000
$ws.Range("N47").Value = -17132
$ws.Range("H105").Value = 7166.3335
$ws.Range("I105").Value = 5750
$ws.Range("K105").Value = 5750
$ws.Range("M105").Value = -4003
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 94
$ws.Range("I12").Value = 185.66667
$ws.Range("J12").Value = 32.88889
$ws.Range("K12").Value = 557.00001
$ws.Range("L12").Value = 98.66667000000001
$ws.Range("M12").Value = -384.00001
$ws.Range("N12").Value = -444.66667
$ws.Range("H39").Value = 781.5
$ws.Range("I39").Value = 781.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2344.5
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2050.5
$ws.Range("N39").Value = ""
$ws.Range("H70").Value = 544
$ws.Range("I70").Value = 544
$ws.Range("K70").Value = 1632
$ws.Range("M70").Value = -1317
$ws.Range("H73").Value = 544
$ws.Range("I73").Value = 544
$ws.Range("K73").Value = 1632
$ws.Range("M73").Value = -540
$ws.Range("H75").Value = 5247.25
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H78").Value = 5247.25
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H112").Value = 7747.5
$ws.Range("I112").Value = 7747.5
$ws.Range("K112").Value = 23242.5
$ws.Range("M112").Value = -22134.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2139
$ws.Range("I107").Value = 706.5714
$ws.Range("K107").Value = 706.5714
$ws.Range("M107").Value = 1213.4286
$ws.Range("H113").Value = 6652.36
$ws.Range("I113").Value = 3750
$ws.Range("J113").Value = 9796.583000000001
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 9796.583000000001
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -14136.583
$ws.Range("H122").Value = 3630.3076
$ws.Range("J122").Value = 4154.375
$ws.Range("L122").Value = 12463.125
$ws.Range("N122").Value = -17363.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2919.7036
$ws.Range("I22").Value = 2305.25
$ws.Range("K22").Value = 2305.25
$ws.Range("M22").Value = -2010.25
$ws.Range("H27").Value = 2919.7036
$ws.Range("I27").Value = 2305.25
$ws.Range("K27").Value = 2305.25
$ws.Range("M27").Value = -2198.25
$ws.Range("H46").Value = 2049.25
$ws.Range("I46").Value = 1399.3334
$ws.Range("J46").Value = 3999
$ws.Range("K46").Value = 1399.3334
$ws.Range("L46").Value = 3999
$ws.Range("M46").Value = -1211.3334
$ws.Range("N46").Value = -4375
$ws.Range("H55").Value = 1645.3334
$ws.Range("I55").Value = 257
$ws.Range("K55").Value = 257
$ws.Range("M55").Value = -84
$ws.Range("H93").Value = 7599.9
$ws.Range("I93").Value = 6000
$ws.Range("J93").Value = 7999.875
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 7999.875
$ws.Range("M93").Value = -4752
$ws.Range("N93").Value = -10495.875
$ws.Range("H100").Value = 7880.1
$ws.Range("I100").Value = 5200.25
$ws.Range("J100").Value = 9666.666999999999
$ws.Range("K100").Value = 5200.25
$ws.Range("L100").Value = 9666.666999999999
$ws.Range("M100").Value = -4659.25
$ws.Range("N100").Value = -10748.667
$ws.Range("H122").Value = 9202.416999999999
$ws.Range("J122").Value = 9500
$ws.Range("L122").Value = 28500
$ws.Range("N122").Value = -33400
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("H107").Value = 3415.5
$ws.Range("I107").Value = 3579.8
$ws.Range("K107").Value = 10739.4
$ws.Range("M107").Value = -8819.400000000001
